# Wisconsin overview workbook: convert numeric "count" cells to text cells
# (values unchanged, but stored/displayed as text, with thousands-separator
# formatting for the "Total" rows) and append a new County "Total" row.
#
# This mirrors the authoring-tool edit described in the commit
# "Update factsheets with text edits from COMM": the No. of 990 Filers
# column (and the Overall-sheet headline count) switch from numeric storage
# to text storage so the values render consistently with the rest of the
# (already-text) sheet, and the County sheet gains a Total summary row.

$wb = $excel.ActiveWorkbook

function Convert-CellToText($cell) {
    # Re-store a numeric cell's value as literal text, preserving the exact
    # displayed digits (with thousands separators for values >= 1000, just
    # like Excel's default General/#,##0 display and like the other totals
    # already present in this workbook).
    $num = $cell.Value2
    $cell.NumberFormat = "@"
    $cell.Value = ("{0:N0}" -f $num)
}

# ---------------------------------------------------------------------
# Sheet "Overall": A2 (headline filer count) 2155 -> "2,155"
# ---------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
Convert-CellToText $wsOverall.Range("A2")

# ---------------------------------------------------------------------
# Sheet "County": B2:B72 per-county counts -> text, plus a new Total row 73
# ---------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")
for ($r = 2; $r -le 72; $r++) {
    Convert-CellToText $wsCounty.Cells.Item($r, 2)
}

$wsCounty.Range("A73").NumberFormat = "@"
$wsCounty.Range("A73").Value = "Total"
$wsCounty.Range("B73").NumberFormat = "@"
$wsCounty.Range("B73").Value = "2,155"
$wsCounty.Range("C73").NumberFormat = "@"
$wsCounty.Range("C73").Value = "`$3,131,962,290"
$wsCounty.Range("D73").NumberFormat = "@"
$wsCounty.Range("D73").Value = "10.44%"
$wsCounty.Range("E73").NumberFormat = "@"
$wsCounty.Range("E73").Value = "-8.25%"
$wsCounty.Range("F73").NumberFormat = "@"
$wsCounty.Range("F73").Value = "63.53%"

# ---------------------------------------------------------------------
# Sheet "Congressional District": B2:B9 counts -> text, B10 Total -> "2,155"
# ---------------------------------------------------------------------
$wsCD = $wb.Worksheets.Item("Congressional District")
for ($r = 2; $r -le 9; $r++) {
    Convert-CellToText $wsCD.Cells.Item($r, 2)
}
Convert-CellToText $wsCD.Cells.Item(10, 2)

# ---------------------------------------------------------------------
# Sheet "Size": B2:B7 counts -> text, B8 Total -> "2,155"
# ---------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")
for ($r = 2; $r -le 7; $r++) {
    Convert-CellToText $wsSize.Cells.Item($r, 2)
}
Convert-CellToText $wsSize.Cells.Item(8, 2)

# ---------------------------------------------------------------------
# Sheet "Subsector": B2:B13 counts -> text, B14 Total -> "2,155"
# ---------------------------------------------------------------------
$wsSubsector = $wb.Worksheets.Item("Subsector")
for ($r = 2; $r -le 13; $r++) {
    Convert-CellToText $wsSubsector.Cells.Item($r, 2)
}
Convert-CellToText $wsSubsector.Cells.Item(14, 2)
